$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (column C) was refreshed for every data row (2-16):
# date serial 46065 -> 46066.
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = 46066
}

# The source data behind rows 7, 8 and 11-16 was re-fetched and the
# records ended up reassigned across those rows. Update "Beteckning"
# (A), "Datum" (B) and "Area (ha)" (G) to match the refreshed data.
$rows = @(7, 8, 11, 12, 13, 14, 15, 16)
$beteckning = @("A 32633-2025", "A 45370-2022", "A 23678-2023", "A 50277-2024", "A 50530-2024", "A 58926-2025", "A 50538-2024", "A 11351-2021")
$datum      = @(45838.65677083333, 44844.6397337963, 45077, 45600.60440972223, 45601.56424768519, 45986, 45601.57153935185, 44263)
$area       = @(1.3, 2.7, 1.4, 0.5, 0.7, 3.1, 0.8, 0.5)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $beteckning[$i]
    $ws.Cells.Item($r, 2).Value = $datum[$i]
    $ws.Cells.Item($r, 7).Value = $area[$i]
}
